$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "adicionado politica de preco"
# Insert two new columns ("modelo" and "politica") before the current
# "full" column (C), shifting full/tipo/link two columns to the right.
$ws.Range("C:D").Insert()

# New column headers
$ws.Range("C1").Value = "modelo"
$ws.Range("D1").Value = "politica"

# "modelo" has no real data yet -> "Sem Modelo" for every product row
$ws.Range("C2:C6").Value = "Sem Modelo"

# "politica" is a new, still-empty column, but the cells must exist as
# (empty) text cells rather than being entirely absent. Using a bare
# quote-prefix forces a text-typed empty cell, then the style is reset
# back to Normal so no quote-prefix formatting is left behind.
$ws.Range("D2:D6").Value = "'"
$ws.Range("D2:D6").Style = "Normal"

# "tipo" values (now column F) are normalized to lowercase
$ws.Range("F2").Value = "premium"
$ws.Range("F3").Value = "premium"
$ws.Range("F4").Value = "premium"
$ws.Range("F5").Value = "classico"
$ws.Range("F6").Value = "premium"

# "link" values (now column G) get refreshed tracking query strings
$ws.Range("G2").Value = "https://produto.mercadolivre.com.br/MLB-2873301438-controle-longa-distancia-jfa-redline-wr-key1-ad1-multimidia-_JM#position%3D13%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Df5bb8ae2-01a8-4ac0-a3e7-9e4810c221df"
$ws.Range("G3").Value = "https://produto.mercadolivre.com.br/MLB-2872406163-controle-longa-distancia-jfa-redline-wr-p-aparelho-original-_JM#position%3D8%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D2694a7a2-cb5d-46bb-b835-b3ff45cf8670"
$ws.Range("G4").Value = "https://produto.mercadolivre.com.br/MLB-2872539303-controle-longa-distncia-k1200-jfa-1200-metros-_JM#position%3D14%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dfcf1f581-3d81-4b93-90fb-631fbd4c42aa"
$ws.Range("G5").Value = "https://produto.mercadolivre.com.br/MLB-3185900332-voltimetro-jfa-vs5hi-3-em-1-sequenciador-high-voltagem-12v-_JM#position%3D8%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3Dc3254e85-464a-47d2-b878-bdca4227f082"
$ws.Range("G6").Value = "https://produto.mercadolivre.com.br/MLB-3185884168-voltimetro-sequenciador-digital-jfa-vs5hi-led-vermelho-_JM#position%3D4%26search_layout%3Dgrid%26type%3Ditem%26tracking_id%3D45721eeb-3c18-4f9f-88fa-49604a41c9c4"
